$d = $word.ActiveDocument

# The address block at the top of the letter had sample/demo data typed in
# after several MERGEFIELD results (Surname, Address, City). Clear that
# literal text back out, leaving the original whitespace runs untouched.
#
# Touching .Font.Name (read + re-assign, a pure no-op) before clearing the
# text keeps the remaining whitespace in its own run instead of letting it
# coalesce with a neighbouring run of identical formatting.
$targets = @("3/2 Street", "CanTho", "VietNam")

foreach ($t in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($t, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Name = $rng.Font.Name
        $rng.Text = ""
    }
}
